# Update "想去人数" (want-to-go count) values in column F across all four
# sheets, matching the upstream data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1173
$ws.Range("F7").Value = 177
$ws.Range("F11").Value = 504
$ws.Range("F14").Value = 726
$ws.Range("F17").Value = 868
$ws.Range("F18").Value = 81495
$ws.Range("F19").Value = 81495
$ws.Range("F22").Value = 35173
$ws.Range("F23").Value = 35173
$ws.Range("F26").Value = 35
$ws.Range("F29").Value = 1071
$ws.Range("F30").Value = 339
$ws.Range("F36").Value = 5576
$ws.Range("F37").Value = 855
$ws.Range("F44").Value = 505
$ws.Range("F45").Value = 11
$ws.Range("F48").Value = 21

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 81
$ws.Range("F35").Value = 27
$ws.Range("F38").Value = 44
$ws.Range("F42").Value = 849
$ws.Range("F43").Value = 315

# Sheet "本地生活" (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 596
$ws.Range("F7").Value = 233

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 596
$ws.Range("F5").Value = 1173
$ws.Range("F11").Value = 177
$ws.Range("F13").Value = 233
$ws.Range("F18").Value = 504
$ws.Range("F21").Value = 726
$ws.Range("F22").Value = 868
$ws.Range("F24").Value = 81495
$ws.Range("F25").Value = 81
$ws.Range("F26").Value = 35173
$ws.Range("F28").Value = 35
$ws.Range("F33").Value = 1071
$ws.Range("F35").Value = 339
$ws.Range("F38").Value = 5576
$ws.Range("F39").Value = 855
$ws.Range("F43").Value = 505
$ws.Range("F44").Value = 11
$ws.Range("F47").Value = 849
$ws.Range("F48").Value = 315
$ws.Range("F49").Value = 21
